# Auto-generated Excel COM-interop script to apply market-data/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Leve" tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 17
$ws.Range("H17").Value = 1856.5714
$ws.Range("J17").Value = 1999.5
$ws.Range("L17").Value = 5998.5
$ws.Range("N17").Value = -6334.5

$ws = $wb.Worksheets.Item("ALC")  # row 64
$ws.Range("H64").Value = 3333.6667
$ws.Range("J64").Value = 3360.6
$ws.Range("L64").Value = 3360.6
$ws.Range("N64").Value = -3856.6

$ws = $wb.Worksheets.Item("ALC")  # row 67
$ws.Range("H67").Value = 3333.6667
$ws.Range("J67").Value = 3360.6
$ws.Range("L67").Value = 3360.6
$ws.Range("N67").Value = -5076.6

$ws = $wb.Worksheets.Item("ALC")  # row 70
$ws.Range("H70").Value = 95233.75
$ws.Range("I70").Value = 1224.5
$ws.Range("K70").Value = 3673.5
$ws.Range("M70").Value = -3403.5

$ws = $wb.Worksheets.Item("ALC")  # row 73
$ws.Range("H73").Value = 95233.75
$ws.Range("I73").Value = 1224.5
$ws.Range("K73").Value = 3673.5
$ws.Range("M73").Value = -2737.5

$ws = $wb.Worksheets.Item("ALC")  # row 106
$ws.Range("H106").Value = 40795.25
$ws.Range("I106").Value = 40795.25
$ws.Range("K106").Value = 40795.25
$ws.Range("M106").Value = -40164.25

$ws = $wb.Worksheets.Item("ALC")  # row 125
$ws.Range("H125").Value = 987.3333
$ws.Range("I125").Value = 985.25
$ws.Range("K125").Value = 8867.25
$ws.Range("M125").Value = -6407.25

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 2783.2856
$ws.Range("I137").Value = 1759.6
$ws.Range("K137").Value = 5278.799999999999
$ws.Range("M137").Value = -2728.799999999999

$ws = $wb.Worksheets.Item("ARM")  # row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 8205.128000000001
$ws.Range("I32").Value = 6176.6763
$ws.Range("J32").Value = 21998.6
$ws.Range("K32").Value = 6176.6763
$ws.Range("L32").Value = 21998.6
$ws.Range("M32").Value = -5889.6763
$ws.Range("N32").Value = -22572.6

$ws = $wb.Worksheets.Item("ARM")  # row 76
$ws.Range("H76").Value = 90000
$ws.Range("J76").Value = 90000
$ws.Range("L76").Value = 90000
$ws.Range("N76").Value = -90676

$ws = $wb.Worksheets.Item("ARM")  # row 79
$ws.Range("H79").Value = 90000
$ws.Range("J79").Value = 90000
$ws.Range("L79").Value = 90000
$ws.Range("N79").Value = -92340

$ws = $wb.Worksheets.Item("BSM")  # row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("BSM")  # row 20
$ws.Range("H20").Value = 771.4
$ws.Range("I20").Value = 771.4
$ws.Range("K20").Value = 771.4
$ws.Range("M20").Value = -524.4

$ws = $wb.Worksheets.Item("BSM")  # row 22
$ws.Range("H22").Value = 632.9167
$ws.Range("I22").Value = 626.8182
$ws.Range("K22").Value = 626.8182
$ws.Range("M22").Value = -453.8182

$ws = $wb.Worksheets.Item("BSM")  # row 26
$ws.Range("H26").Value = 22499.5
$ws.Range("I26").Value = 22499.5
$ws.Range("K26").Value = 22499.5
$ws.Range("M26").Value = -22207.5

$ws = $wb.Worksheets.Item("BSM")  # row 42
$ws.Range("H42").Value = 170000
$ws.Range("J42").Value = 170000
$ws.Range("L42").Value = 170000
$ws.Range("N42").Value = -170656

$ws = $wb.Worksheets.Item("BSM")  # row 64
$ws.Range("H64").Value = 697.58826
$ws.Range("I64").Value = 567.4286
$ws.Range("J64").Value = 788.7
$ws.Range("K64").Value = 567.4286
$ws.Range("L64").Value = 788.7
$ws.Range("M64").Value = -342.4286
$ws.Range("N64").Value = -1238.7

$ws = $wb.Worksheets.Item("BSM")  # row 67
$ws.Range("H67").Value = 697.58826
$ws.Range("I67").Value = 567.4286
$ws.Range("J67").Value = 788.7
$ws.Range("K67").Value = 567.4286
$ws.Range("L67").Value = 788.7
$ws.Range("M67").Value = 212.5714
$ws.Range("N67").Value = -2348.7

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 1627.7333
$ws.Range("I134").Value = 1627.7333
$ws.Range("K134").Value = 4883.199900000001
$ws.Range("M134").Value = -2348.199900000001

$ws = $wb.Worksheets.Item("CRP")  # row 31
$ws.Range("H31").Value = 6296.533
$ws.Range("I31").Value = 4151
$ws.Range("K31").Value = 4151
$ws.Range("M31").Value = -3856

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 6296.533
$ws.Range("I34").Value = 4151
$ws.Range("K34").Value = 4151
$ws.Range("M34").Value = -3949

$ws = $wb.Worksheets.Item("CRP")  # row 69
$ws.Range("H69").Value = 19428
$ws.Range("I69").Value = 19428
$ws.Range("K69").Value = 19428
$ws.Range("M69").Value = -18679

$ws = $wb.Worksheets.Item("CRP")  # row 72
$ws.Range("H72").Value = 19428
$ws.Range("I72").Value = 19428
$ws.Range("K72").Value = 58284
$ws.Range("M72").Value = -54540

$ws = $wb.Worksheets.Item("CRP")  # row 74
$ws.Range("H74").Value = 47782
$ws.Range("J74").Value = 45314
$ws.Range("L74").Value = 45314
$ws.Range("N74").Value = -47062

$ws = $wb.Worksheets.Item("CRP")  # row 77
$ws.Range("H77").Value = 47782
$ws.Range("J77").Value = 45314
$ws.Range("L77").Value = 135942
$ws.Range("N77").Value = -144678

$ws = $wb.Worksheets.Item("CRP")  # row 92
$ws.Range("H92").Value = 23900.25
$ws.Range("J92").Value = 23900.25
$ws.Range("L92").Value = 23900.25
$ws.Range("N92").Value = -28892.25

$ws = $wb.Worksheets.Item("CRP")  # row 103
$ws.Range("H103").Value = 11809.4
$ws.Range("I103").Value = 13761.75
$ws.Range("K103").Value = 13761.75
$ws.Range("M103").Value = -12589.75

$ws = $wb.Worksheets.Item("CRP")  # row 122
$ws.Range("H122").Value = 5060.75
$ws.Range("I122").Value = 4898.1333
$ws.Range("K122").Value = 14694.3999
$ws.Range("M122").Value = -12244.3999

$ws = $wb.Worksheets.Item("CRP")  # row 132
$ws.Range("H132").Value = 3916.375
$ws.Range("I132").Value = 1949.75
$ws.Range("J132").Value = 5883
$ws.Range("K132").Value = 5849.25
$ws.Range("L132").Value = 17649
$ws.Range("M132").Value = -3319.25
$ws.Range("N132").Value = -22709

$ws = $wb.Worksheets.Item("CUL")  # row 128
$ws.Range("H128").Value = 3979899
$ws.Range("I128").Value = 3979899
$ws.Range("K128").Value = 11939697
$ws.Range("M128").Value = -11934717

$ws = $wb.Worksheets.Item("GSM")  # row 70
$ws.Range("H70").Value = 5427.7144
$ws.Range("J70").Value = 4999.2
$ws.Range("L70").Value = 4999.2
$ws.Range("N70").Value = -5539.2

$ws = $wb.Worksheets.Item("GSM")  # row 73
$ws.Range("H73").Value = 5427.7144
$ws.Range("J73").Value = 4999.2
$ws.Range("L73").Value = 4999.2
$ws.Range("N73").Value = -6871.2

$ws = $wb.Worksheets.Item("GSM")  # row 113
$ws.Range("H113").Value = 3784.3333
$ws.Range("I113").Value = 2123.6
$ws.Range("J113").Value = 4970.5713
$ws.Range("K113").Value = 2123.6
$ws.Range("L113").Value = 4970.5713
$ws.Range("M113").Value = 46.40000000000009
$ws.Range("N113").Value = -9310.5713

$ws = $wb.Worksheets.Item("LTW")  # row 6
$ws.Range("H6").Value = 50000
$ws.Range("J6").Value = 50000
$ws.Range("L6").Value = 50000
$ws.Range("N6").Value = -50224

$ws = $wb.Worksheets.Item("LTW")  # row 46
$ws.Range("H46").Value = 1268.8
$ws.Range("I46").Value = 1268.8
$ws.Range("K46").Value = 1268.8
$ws.Range("M46").Value = -1080.8

$ws = $wb.Worksheets.Item("LTW")  # row 82
$ws.Range("H82").Value = 6481
$ws.Range("J82").Value = 5999.6
$ws.Range("L82").Value = 5999.6
$ws.Range("N82").Value = -6721.6

$ws = $wb.Worksheets.Item("LTW")  # row 85
$ws.Range("H85").Value = 6481
$ws.Range("J85").Value = 5999.6
$ws.Range("L85").Value = 5999.6
$ws.Range("N85").Value = -8495.6

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Range("H122").Value = 10329.8
$ws.Range("I122").Value = 9162.25
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 27486.75
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -25036.75
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("LTW")  # row 136
$ws.Range("H136").Value = 2500
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")  # row 113
$ws.Range("H113").Value = 584.8570999999999
$ws.Range("I113").Value = 781
$ws.Range("J113").Value = 231.8
$ws.Range("K113").Value = 2343
$ws.Range("L113").Value = 695.4000000000001
$ws.Range("M113").Value = -173
$ws.Range("N113").Value = -5035.4

$ws = $wb.Worksheets.Item("WVR")  # row 122
$ws.Range("H122").Value = 8340191
$ws.Range("I122").Value = 8340191
$ws.Range("K122").Value = 25020573
$ws.Range("M122").Value = -25018123

$ws = $wb.Worksheets.Item("WVR")  # row 132
$ws.Range("H132").Value = 2232.4285
$ws.Range("I132").Value = 1477.7273
$ws.Range("K132").Value = 4433.1819
$ws.Range("M132").Value = -1903.1819

